$d = $word.ActiveDocument

# "Full Stack .NET Developer with 1.9 years of experience..."
# becomes
# "Full Stack .NET Developer with 2 years of experience..."
$d.Content.Find.Execute("1.9 years", $true, $false, $false, $false, $false, $true, 1, $false, "2 years", 2)
